$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 494, pushing existing rows 494-615 down to 495-616.
$ws.Rows.Item(494).Insert()

# Populate the newly inserted row 494 with its data (mirrors the sibling rows,
# with new Fecha/Volumen/Precio values per the commit).
$ws.Range("A494").Value = 3
$ws.Range("B494").Value = "Femacal de La Calera"
$ws.Range("C494").Value = "Coquimbo"
$ws.Range("D494").Value = 44932
$ws.Range("E494").Value = 5
$ws.Range("F494").Value = 100112037
$ws.Range("G494").Value = "Cebollín"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 170
$ws.Range("K494").Value = 3500
$ws.Range("L494").Value = 4000
$ws.Range("M494").Value = 3765
$ws.Range("N494").Value = "$/paquete 36 unidades"
$ws.Range("O494").Value = "Provincia de Quillota"
$ws.Range("P494").Value = 105
$ws.Range("Q494").Value = 36
$ws.Range("R494").Value = "Hortaliza"
